# ModelRuns_RTP2025.xlsx — add 2023 v41 and v42 runs to the run log
#
# Inserts two new rows (54 and 55) into the "ModelRuns" sheet, just
# before the existing 2025 row ("25_TM152_FBP_Plus_22"), for the new
# model runs 2023_TM160_IPA_41 and 2023_TM160_IPA_42. This pushes the
# previously-existing rows 54-64 down to rows 56-66.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 54, inheriting that
# row's formatting (same as Excel's native "Insert Sheet Rows").
$ws.Rows.Item(54).Resize(2).Insert()

# --- Row 54: 2023_TM160_IPA_41 -------------------------------------
$ws.Range("A54").Value = 2023
$ws.Range("B54").Value = "2023_TM160_IPA_41"
$ws.Range("C54").Value = "RTP2025_IP"
$ws.Range("D54").Value = "Past year"
$ws.Range("E54").Value = "AOC=16.61, WFH at ~30%"
$ws.Range("F54").Value = "petrale"
$ws.Range("G54").Value = "n/a"
$ws.Range("H54").Value = "current"
$ws.Range("I54").Value = "BlueprintNetworks_v12\net_2023_Blueprint"
$ws.Range("J54").Value = "model3-c"
$ws.Range("K54").Value = "https://app.asana.com/0/1204085012544660/1205980528918816/f"
$ws.Range("L54").Value = 16.61
$ws.Range("M54").Value = "na"
$ws.Range("N54").Value = "na"
$ws.Range("O54").Value = 0.99
$ws.Range("P54").Value = 0.89
$ws.Range("Q54").Value = 100
$ws.Range("R54").Value = 0
$ws.Range("S54").Value = 75

# --- Row 55: 2023_TM160_IPA_42 -------------------------------------
$ws.Range("A55").Value = 2023
$ws.Range("B55").Value = "2023_TM160_IPA_42"
$ws.Range("C55").Value = "RTP2025_IP"
$ws.Range("D55").Value = "Past year"
$ws.Range("E55").Value = "AOC=16.61, WFH at ~31%"
$ws.Range("F55").Value = "petrale"
$ws.Range("G55").Value = "n/a"
$ws.Range("H55").Value = "current"
$ws.Range("I55").Value = "BlueprintNetworks_v12\net_2023_Blueprint"
$ws.Range("J55").Value = "model3-b"
$ws.Range("K55").Value = "https://app.asana.com/0/1204085012544660/1205983427401938/f"
$ws.Range("L55").Value = 16.61
$ws.Range("M55").Value = "na"
$ws.Range("N55").Value = "na"
$ws.Range("O55").Value = 1.04
$ws.Range("P55").Value = 0.94
$ws.Range("Q55").Value = 100
$ws.Range("R55").Value = 0
$ws.Range("S55").Value = 75

# --- Misc cleanup ----------------------------------------------------
# K49 ("2023_TM160_IPA_40" Asana link) previously carried a redundant
# "applyFill" style distinct from the sheet's plain text style; reset
# its fill so it matches the plain style used by the rest of column K.
$ws.Range("K49").Interior.Pattern = -4142  # xlPatternNone

# Update the active selection to mirror where the editor's cursor
# ended up after the edit.
$ws.Range("J56").Select() | Out-Null
